$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trait rows: related observation load functions for LAB color + fruit weight
$rows = @(
    @("LAB Color1", "LAB_color"),
    @("LAB Color2", "LAB_color"),
    @("LAB Color3", "LAB_color"),
    @("LAB Color4", "LAB_color"),
    @("number_fruit", "numeric"),
    @("verage_Fruit_weight_g", "numeric"),
    @("Average_Fruit_weight_g", "numeric"),
    @("Fruit_weight_g", "numeric")
)

$r = 5
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

$null = $ws.Range("A9:A12").Select()
